$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates - force text storage so formats like trailing zeros / multi-dot strings are preserved
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.563.54'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.883.71'
$ws.Range("D3").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4740'
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06540'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.31'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7749'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '100.93'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07814'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.883.93'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.258'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '284.66'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.554.38'
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007527'
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.130.14'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.359'
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.440'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.169'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.89'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.12'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.915'
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09704'
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.258'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.190'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04845'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.130'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6981'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.757'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01913'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.900'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '76.02'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.297'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.985'
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8306'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.69'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.885'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.019'
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '894.51'
$ws.Range("D50").Style = "Normal"

# Volume(1h) (column E) updates - plain text percentage strings
$ws.Range("E2").Value = '  -0.56%  '
$ws.Range("E3").Value = '  -0.30%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("E5").Value = '  -0.66%  '
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("E8").Value = '  -1.05%  '
$ws.Range("E9").Value = '  +0.08%  '
$ws.Range("E10").Value = '  +0.86%  '
$ws.Range("E11").Value = '  +4.93%  '
$ws.Range("E12").Value = '  +4.17%  '
$ws.Range("E13").Value = '  +0.13%  '
$ws.Range("E14").Value = '  -0.24%  '
$ws.Range("E15").Value = '  +0.19%  '
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("E17").Value = '  -0.53%  '
$ws.Range("E18").Value = '  -0.39%  '
$ws.Range("E19").Value = '  -0.31%  '
$ws.Range("E20").Value = '  +0.02%  '
$ws.Range("E21").Value = '  -0.27%  '
$ws.Range("E22").Value = '  +0.62%  '
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("E24").Value = '  +2.95%  '
$ws.Range("E25").Value = '  -0.67%  '
$ws.Range("E26").Value = '  -1.25%  '
$ws.Range("E27").Value = '  +0.57%  '
$ws.Range("E29").Value = '  -0.41%  '
$ws.Range("E30").Value = '  -0.48%  '
$ws.Range("E31").Value = '  +0.75%  '
$ws.Range("E32").Value = '  -0.98%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("E34").Value = '  -0.37%  '
$ws.Range("E35").Value = '  +0.35%  '
$ws.Range("E36").Value = '  +0.07%  '
$ws.Range("E37").Value = '  +1.25%  '
$ws.Range("E38").Value = '  +0.97%  '
$ws.Range("E39").Value = '  +3.37%  '
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("E41").Value = '  -0.46%  '
$ws.Range("E42").Value = '  -0.74%  '
$ws.Range("E43").Value = '  -0.50%  '
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("E45").Value = '  -0.62%  '
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("E47").Value = '  +3.69%  '
$ws.Range("E48").Value = '  -0.55%  '
$ws.Range("E49").Value = '  -1.19%  '
$ws.Range("E50").Value = '  -2.74%  '
